{"js": "// Remove the unused \"Abstract Title\" paragraph style and bump the\n// \"Abstract\" style's space-before to match its space-after (15pt == 300\n// twentieths-of-a-point), per the target OOXML diff.\n\nconst styles = context.document.getStyles();\n\n// 1) Delete the \"AbstractTitle\" custom style entirely.\nconst abstractTitleStyle = styles.getByNameOrNullObject(\"Abstract Title\");\nabstractTitleStyle.load(\"isNullObject\");\nawait context.sync();\n\nif (!abstractTitleStyle.isNullObject) {\n  abstractTitleStyle.delete();\n}\n\n// 2) Change \"Abstract\" style spacing: before 100 -> 300 (twentieths of a\n// point), i.e. 5pt -> 15pt in the Office.js (points-based) API. Leave\n// spaceAfter (300 twentieths == 15pt) untouched.\nconst abstractStyle = styles.getByNameOrNullObject(\"Abstract\");\nabstractStyle.load(\"isNullObject\");\nawait context.sync();\n\nif (!abstractStyle.isNullObject) {\n  abstractStyle.paragraphFormat.spaceBefore = 15;\n}\n\nawait context.sync();\n", "ps1": "# Remove the unused \"Abstract Title\" paragraph style and bump the\n# \"Abstract\" style's space-before to match its space-after (15pt == 300\n# twentieths-of-a-point), per the target OOXML diff.\n\n$d = $word.ActiveDocument\n\n# 1) Delete the \"AbstractTitle\" custom style entirely.\n$titleStyle = $d.Styles(\"Abstract Title\")\nif ($titleStyle -ne $null) {\n    $titleStyle.Delete()\n}\n\n# 2) Change \"Abstract\" style spacing: before 100 -> 300 (twentieths of a\n# point), i.e. 5pt -> 15pt in the Word COM (points-based) API. Leave\n# SpaceAfter (300 twentieths == 15pt) untouched.\n$abstractStyle = $d.Styles(\"Abstract\")\n$abstractStyle.ParagraphFormat.SpaceBefore = 15\n"}
